{"js": "// Replace the date line and each \"A\u00f7B=C, D\" cell value with its updated\n// value. Every old string is unique within the document, so a simple\n// exact-text search + replace (one hit each) is unambiguous and safe.\nconst replacements = [\n  [\"2025-07-02 Wednesday\", \"2025-07-03 Thursday\"],\n  [\"514\u00f76=85, 4\", \"808\u00f79=89, 7\"],\n  [\"474\u00f73=158, 0\", \"777\u00f74=194, 1\"],\n  [\"666\u00f72=333, 0\", \"948\u00f74=237, 0\"],\n  [\"386\u00f73=128, 2\", \"110\u00f78=13, 6\"],\n  [\"671\u00f76=111, 5\", \"108\u00f79=12, 0\"],\n  [\"769\u00f74=192, 1\", \"876\u00f75=175, 1\"],\n  [\"749\u00f77=107, 0\", \"113\u00f77=16, 1\"],\n  [\"738\u00f72=369, 0\", \"314\u00f79=34, 8\"],\n  [\"459\u00f79=51, 0\", \"981\u00f79=109, 0\"],\n  [\"522\u00f73=174, 0\", \"551\u00f76=91, 5\"],\n  [\"821\u00f79=91, 2\", \"921\u00f79=102, 3\"],\n  [\"415\u00f73=138, 1\", \"273\u00f77=39, 0\"],\n  [\"483\u00f74=120, 3\", \"218\u00f74=54, 2\"],\n  [\"269\u00f77=38, 3\", \"552\u00f77=78, 6\"],\n  [\"298\u00f77=42, 4\", \"307\u00f74=76, 3\"],\n  [\"540\u00f76=90, 0\", \"460\u00f78=57, 4\"],\n  [\"523\u00f79=58, 1\", \"615\u00f72=307, 1\"],\n  [\"629\u00f77=89, 6\", \"931\u00f77=133, 0\"],\n  [\"782\u00f72=391, 0\", \"319\u00f76=53, 1\"],\n  [\"235\u00f78=29, 3\", \"191\u00f75=38, 1\"],\n  [\"757\u00f79=84, 1\", \"378\u00f78=47, 2\"],\n  [\"715\u00f75=143, 0\", \"124\u00f76=20, 4\"],\n  [\"245\u00f77=35, 0\", \"816\u00f75=163, 1\"],\n  [\"373\u00f74=93, 1\", \"426\u00f77=60, 6\"],\n  [\"857\u00f73=285, 2\", \"622\u00f78=77, 6\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the header date and each \"A\u00f7B=C, D\" table cell to its new value.\n# Every old string is unique in the document, so Find/Replace (one hit\n# each, ReplaceAll) is unambiguous and safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-07-02 Wednesday\", \"2025-07-03 Thursday\"),\n    @(\"514\u00f76=85, 4\", \"808\u00f79=89, 7\"),\n    @(\"474\u00f73=158, 0\", \"777\u00f74=194, 1\"),\n    @(\"666\u00f72=333, 0\", \"948\u00f74=237, 0\"),\n    @(\"386\u00f73=128, 2\", \"110\u00f78=13, 6\"),\n    @(\"671\u00f76=111, 5\", \"108\u00f79=12, 0\"),\n    @(\"769\u00f74=192, 1\", \"876\u00f75=175, 1\"),\n    @(\"749\u00f77=107, 0\", \"113\u00f77=16, 1\"),\n    @(\"738\u00f72=369, 0\", \"314\u00f79=34, 8\"),\n    @(\"459\u00f79=51, 0\", \"981\u00f79=109, 0\"),\n    @(\"522\u00f73=174, 0\", \"551\u00f76=91, 5\"),\n    @(\"821\u00f79=91, 2\", \"921\u00f79=102, 3\"),\n    @(\"415\u00f73=138, 1\", \"273\u00f77=39, 0\"),\n    @(\"483\u00f74=120, 3\", \"218\u00f74=54, 2\"),\n    @(\"269\u00f77=38, 3\", \"552\u00f77=78, 6\"),\n    @(\"298\u00f77=42, 4\", \"307\u00f74=76, 3\"),\n    @(\"540\u00f76=90, 0\", \"460\u00f78=57, 4\"),\n    @(\"523\u00f79=58, 1\", \"615\u00f72=307, 1\"),\n    @(\"629\u00f77=89, 6\", \"931\u00f77=133, 0\"),\n    @(\"782\u00f72=391, 0\", \"319\u00f76=53, 1\"),\n    @(\"235\u00f78=29, 3\", \"191\u00f75=38, 1\"),\n    @(\"757\u00f79=84, 1\", \"378\u00f78=47, 2\"),\n    @(\"715\u00f75=143, 0\", \"124\u00f76=20, 4\"),\n    @(\"245\u00f77=35, 0\", \"816\u00f75=163, 1\"),\n    @(\"373\u00f74=93, 1\", \"426\u00f77=60, 6\"),\n    @(\"857\u00f73=285, 2\", \"622\u00f78=77, 6\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
